# Add a new "2022-Q1" sheet (positioned after "2021-Q2", before "总计")
# and record its summary stats on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$q2sheet = $wb.Worksheets.Item("2021-Q2")
# A default-styled cell (no explicit style) used below to "reset" style
# indexes back to the workbook default after text-forcing a cell via
# NumberFormat.
$plainCell = $q2sheet.Range("C2")

# --- 1. Create the new "2022-Q1" sheet right after "2021-Q2" ------------
$newSheet = $wb.Worksheets.Add($null, $q2sheet)
$newSheet.Name = "2022-Q1"

# Header row (B1:H1) -- same labels/layout as the other quarterly sheets
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the bold/bordered header style from the "2021-Q2" sheet
$q2sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Data rows (2-6). Columns B-G are stored as text (matches the other
# quarterly sheets, and keeps fund codes like "000369"/"003718" from
# losing their leading zeros). Force text via NumberFormat, then strip
# the resulting custom number-format style back off with a format-only
# paste from a plain (default-style) cell.
$textRange = $newSheet.Range("B2:G6")
$textRange.NumberFormat = "@"

$newSheet.Range("B2").Value = "513500"
$newSheet.Range("C2").Value = "博时标普500ETF(QDII)"
$newSheet.Range("D2").Value = "70.03"
$newSheet.Range("E2").Value = "90.45"
$newSheet.Range("F2").Value = "1.13"
$newSheet.Range("G2").Value = "0.7913"
$newSheet.Range("H2").Value = 10

$newSheet.Range("B3").Value = "000369"
$newSheet.Range("C3").Value = "广发全球医疗保健(QDII) - 人民币"
$newSheet.Range("D3").Value = "2.46"
$newSheet.Range("E3").Value = "81.85"
$newSheet.Range("F3").Value = "5.79"
$newSheet.Range("G3").Value = "0.1424"
$newSheet.Range("H3").Value = 2

$newSheet.Range("B4").Value = "000370"
$newSheet.Range("C4").Value = "广发全球医疗保健(QDII) - 美元"
$newSheet.Range("D4").Value = "2.46"
$newSheet.Range("E4").Value = "81.85"
$newSheet.Range("F4").Value = "5.79"
$newSheet.Range("G4").Value = "0.1424"
$newSheet.Range("H4").Value = 2

$newSheet.Range("B5").Value = "003718"
$newSheet.Range("C5").Value = "易方达标普500指数(QDII-LOF) 美元"
$newSheet.Range("D5").Value = "5.22"
$newSheet.Range("E5").Value = "91.11"
$newSheet.Range("F5").Value = "1.14"
$newSheet.Range("G5").Value = "0.0595"
$newSheet.Range("H5").Value = 10

$newSheet.Range("B6").Value = "161125"
$newSheet.Range("C6").Value = "易方达标普500指数(QDII-LOF) 人民币"
$newSheet.Range("D6").Value = "5.22"
$newSheet.Range("E6").Value = "91.11"
$newSheet.Range("F6").Value = "1.14"
$newSheet.Range("G6").Value = "0.0595"
$newSheet.Range("H6").Value = 10

$plainCell.Copy()
$textRange.PasteSpecial(-4122)

# Column A index (0,1,2,3,4) with the same bold/bordered style used
# on the other quarterly sheets
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4

$q2sheet.Range("A2:A6").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

# --- 2. Update the "总计" (totals) sheet with the new 2022-Q1 row -------
$total = $wb.Worksheets.Item("总计")

$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.2

# Renumber the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the bold/bordered index-column style on the new row
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
